$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 19.666666
$ws.Range("I42").Value = 19.666666
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 58.999998
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 171.000002
$ws.Range("N42").ClearContents()
$ws.Range("H113").Value = 1764.091
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 1711.6666
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 1711.6666
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -8219.6666
$ws.Range("H137").Value = 29413714
$ws.Range("I137").Value = 1327.8846
$ws.Range("J137").Value = 125003976
$ws.Range("K137").Value = 3983.6538
$ws.Range("L137").Value = 375011928
$ws.Range("M137").Value = -1433.6538
$ws.Range("N137").Value = -375017028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2652.3225
$ws.Range("I61").Value = 1807.76
$ws.Range("J61").Value = 6171.3335
$ws.Range("K61").Value = 1807.76
$ws.Range("L61").Value = 6171.3335
$ws.Range("M61").Value = -1595.76
$ws.Range("N61").Value = -6595.3335
$ws.Range("H74").Value = 5980.8
$ws.Range("I74").Value = 848.26666
$ws.Range("J74").Value = 11113.333
$ws.Range("K74").Value = 848.26666
$ws.Range("L74").Value = 11113.333
$ws.Range("M74").Value = 25.73334
$ws.Range("N74").Value = -12861.333
$ws.Range("H77").Value = 5980.8
$ws.Range("I77").Value = 848.26666
$ws.Range("J77").Value = 11113.333
$ws.Range("K77").Value = 4241.3333
$ws.Range("L77").Value = 55566.665
$ws.Range("M77").Value = 126.6666999999998
$ws.Range("N77").Value = -64302.665
$ws.Range("H132").Value = 3325.3684
$ws.Range("I132").Value = 2293.0908
$ws.Range("J132").Value = 4744.75
$ws.Range("K132").Value = 6879.2724
$ws.Range("L132").Value = 14234.25
$ws.Range("M132").Value = -4349.2724
$ws.Range("N132").Value = -19294.25
$ws.Range("H136").Value = 2652.3225
$ws.Range("I136").Value = 1807.76
$ws.Range("J136").Value = 6171.3335
$ws.Range("K136").Value = 5423.28
$ws.Range("L136").Value = 18514.0005
$ws.Range("M136").Value = -2873.28
$ws.Range("N136").Value = -23614.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 37661
$ws.Range("I134").Value = 50125.76
$ws.Range("J134").Value = 4941
$ws.Range("K134").Value = 150377.28
$ws.Range("L134").Value = 14823
$ws.Range("M134").Value = -147842.28
$ws.Range("N134").Value = -19893

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1414.8636
$ws.Range("I31").Value = 1306.6842
$ws.Range("J31").Value = 2100
$ws.Range("K31").Value = 1306.6842
$ws.Range("L31").Value = 2100
$ws.Range("M31").Value = -1011.6842
$ws.Range("N31").Value = -2690
$ws.Range("H34").Value = 1414.8636
$ws.Range("I34").Value = 1306.6842
$ws.Range("J34").Value = 2100
$ws.Range("K34").Value = 1306.6842
$ws.Range("L34").Value = 2100
$ws.Range("M34").Value = -1104.6842
$ws.Range("N34").Value = -2504
$ws.Range("H58").Value = 2656.4583
$ws.Range("I58").Value = 1360.9231
$ws.Range("J58").Value = 4187.5454
$ws.Range("K58").Value = 1360.9231
$ws.Range("L58").Value = 4187.5454
$ws.Range("M58").Value = -1157.9231
$ws.Range("N58").Value = -4593.5454
$ws.Range("H94").Value = 1146.5385
$ws.Range("I94").Value = 1049.25
$ws.Range("J94").Value = 1189.7778
$ws.Range("K94").Value = 1049.25
$ws.Range("L94").Value = 1189.7778
$ws.Range("M94").Value = -598.25
$ws.Range("N94").Value = -2091.7778
$ws.Range("H132").Value = 2419.7727
$ws.Range("I132").Value = 1575.4546
$ws.Range("J132").Value = 3264.0908
$ws.Range("K132").Value = 4726.3638
$ws.Range("L132").Value = 9792.2724
$ws.Range("M132").Value = -2196.3638
$ws.Range("N132").Value = -14852.2724
$ws.Range("H134").Value = 2002.55
$ws.Range("I134").Value = 2252.0688
$ws.Range("J134").Value = 1344.7273
$ws.Range("K134").Value = 6756.2064
$ws.Range("L134").Value = 4034.1819
$ws.Range("M134").Value = -4221.2064
$ws.Range("N134").Value = -9104.1819
$ws.Range("H136").Value = 2656.4583
$ws.Range("I136").Value = 1360.9231
$ws.Range("J136").Value = 4187.5454
$ws.Range("K136").Value = 4082.7693
$ws.Range("L136").Value = 12562.6362
$ws.Range("M136").Value = -1532.7693
$ws.Range("N136").Value = -17662.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 69.80769
$ws.Range("I12").Value = 99.5
$ws.Range("J12").Value = 44.357143
$ws.Range("K12").Value = 298.5
$ws.Range("L12").Value = 133.071429
$ws.Range("M12").Value = -125.5
$ws.Range("N12").Value = -479.071429
$ws.Range("H94").Value = 3100.4
$ws.Range("I94").Value = 1274
$ws.Range("J94").Value = 4318
$ws.Range("K94").Value = 3822
$ws.Range("L94").Value = 12954
$ws.Range("M94").Value = -3146
$ws.Range("N94").Value = -14306
$ws.Range("H100").Value = 2768.2856
$ws.Range("J100").Value = 2899.6667
$ws.Range("L100").Value = 8699.000100000001
$ws.Range("N100").Value = -10321.0001
$ws.Range("H105").Value = 181602400
$ws.Range("J105").Value = 181602400
$ws.Range("L105").Value = 544807200
$ws.Range("N105").Value = -544812442
$ws.Range("H129").Value = 919.5789
$ws.Range("J129").Value = 1089.4286
$ws.Range("L129").Value = 3268.2858
$ws.Range("N129").Value = -13268.2858
$ws.Range("H137").Value = 26750.512
$ws.Range("I137").Value = 2472.3076
$ws.Range("J137").Value = 36613.53
$ws.Range("K137").Value = 7416.9228
$ws.Range("L137").Value = 109840.59
$ws.Range("M137").Value = -2316.9228
$ws.Range("N137").Value = -120040.59

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2008.0714
$ws.Range("I132").Value = 1212.625
$ws.Range("J132").Value = 3068.6667
$ws.Range("K132").Value = 3637.875
$ws.Range("L132").Value = 9206.000100000001
$ws.Range("M132").Value = -1107.875
$ws.Range("N132").Value = -14266.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7008.522
$ws.Range("I132").Value = 8262.375
$ws.Range("J132").Value = 4142.5713
$ws.Range("K132").Value = 24787.125
$ws.Range("L132").Value = 12427.7139
$ws.Range("M132").Value = -22257.125
$ws.Range("N132").Value = -17487.7139
$ws.Range("H136").Value = 2536.4707
$ws.Range("I136").Value = 1325.5555
$ws.Range("J136").Value = 3898.75
$ws.Range("K136").Value = 3976.6665
$ws.Range("L136").Value = 11696.25
$ws.Range("M136").Value = -1426.6665
$ws.Range("N136").Value = -16796.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3506.5715
$ws.Range("I132").Value = 4893.067
$ws.Range("J132").Value = 1906.7693
$ws.Range("K132").Value = 14679.201
$ws.Range("L132").Value = 5720.3079
$ws.Range("M132").Value = -12149.201
$ws.Range("N132").Value = -10780.3079
$ws.Range("H136").Value = 8421.25
$ws.Range("I136").Value = 10322.728
$ws.Range("J136").Value = 1449.1666
$ws.Range("K136").Value = 30968.184
$ws.Range("L136").Value = 4347.4998
$ws.Range("M136").Value = -28418.184
$ws.Range("N136").Value = -9447.4998
